$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '97.032.74'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.23%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.687.53'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.06%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.44%  '

# Row 6
$ws.Range("E6").Value = '  +9.74%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '655.18'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.08%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.425'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.57%  '

# Row 9
$ws.Range("E9").Value = '  +3.77%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.999'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.00%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.685.92'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.09%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.61'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.06%  '

# Row 13
$ws.Range("E13").Value = '  +1.29%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.91'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.58%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.369.31'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.98%  '

# Row 16
$ws.Range("E16").Value = '  +3.27%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '96.770.35'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.20%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.08'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.58%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.676.23'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.94%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.10%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.80'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.83%  '

# Row 22
$ws.Range("E22").Value = '  -0.27%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '532.22'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.41%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.51'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.79%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.16'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.26%  '

# Row 26
$ws.Range("E26").Value = '  -0.62%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '102.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.77%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '13.52'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.50%  '

# Row 29
$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.168'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.12%  '

# Row 30
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.52'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.47%  '

# Row 31
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.05'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.75%  '

# Row 32
$ws.Range("B32").Value = 'Dai'
$ws.Range("C32").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.02%  '

# Row 33
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.90'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +15.75%  '

# Row 34
$ws.Range("B34").Value = 'Cronos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.186'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.45%  '

# Row 35
$ws.Range("B35").Value = 'Binance-PegBSC-USD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.01'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.95%  '

# Row 36
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '32.80'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.82%  '

# Row 37
$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '655.86'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.62%  '

# Row 38
$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.607'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.43%  '

# Row 39
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.04'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.34%  '

# Row 40
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.96'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +15.36%  '

# Row 41
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.164'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.47%  '

# Row 42
$ws.Range("B42").Value = 'ImmutableX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.16%  '

# Row 43
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.967'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.56%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '38.43'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +16.73%  '

# Row 45
$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.02%  '

# Row 46
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.453'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.85%  '

# Row 47
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0461'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.79%  '

# Row 48
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.34'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.62%  '

# Row 49
$ws.Range("B49").Value = 'WhiteBITCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.64'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.12%  '

# Row 50
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.76'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.57%  '

# Row 51
$ws.Range("B51").Value = 'MantraDAO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.64'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.91%  '
